# Updated symbol list (Price / Volume(1h) columns) to reflect latest
# crypto-ranking snapshot.
#
# All of these cells hold plain text (not numbers/percent-formatted
# numbers) in the workbook, so each value is written with a leading
# apostrophe to force text entry (preventing Excel from re-interpreting
# e.g. "0.1390" as the number 0.139 and dropping the trailing zero, or
# "1.16%" as a percentage value). ClearFormats() immediately afterwards
# strips the "Text" number-format Excel auto-applies when you force text
# entry like that, so the cell's style stays exactly as it was (no
# explicit number format), matching the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.31"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'1.16%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'35.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'1.28%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.054"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.14%"
$ws.Range("E4").ClearFormats()
$ws.Range("E5").Value = "'0.89%"
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'0.74%"
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'3.14%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'7.725"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-1.03%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.9276"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'0.60%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1390"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'9.00%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.1895"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'2.53%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.09228"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-5.15%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.03598"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'0.86%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.09812"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.40%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.001411"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'1.41%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.005901"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'2.35%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.552"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'1.35%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.880"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-1.73%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.3468"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "'0.1303"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'0.87%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'-3.11%"
$ws.Range("E21").ClearFormats()
$ws.Range("D23").Value = "'0.04436"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-1.54%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001221"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'0.75%"
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'0.03%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.0001651"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'31.85%"
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.0003129"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'4.23%"
$ws.Range("E27").ClearFormats()
$ws.Range("D39").Value = "'0.01955"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'4.34%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.04917"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'4.76%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007630"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'2.40%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'-8.35%"
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'3.67%"
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'-0.57%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.01139"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'6.50%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006381"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'2.11%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'65.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'1.15%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'-20.05%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'-0.12%"
$ws.Range("E51").ClearFormats()
